$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2657
$ws1.Range("F3").Value = 339
$ws1.Range("F4").Value = 353
$ws1.Range("F5").Value = 1485
$ws1.Range("F6").Value = 1134
$ws1.Range("F7").Value = 326
$ws1.Range("F11").Value = 115
$ws1.Range("F13").Value = 9071
$ws1.Range("F14").Value = 391
$ws1.Range("F15").Value = 2500
$ws1.Range("F16").Value = 2
$ws1.Range("F17").Value = 255
$ws1.Range("F20").Value = 623
$ws1.Range("F22").Value = 1177
$ws1.Range("F24").Value = 2081
$ws1.Range("F25").Value = 2165
$ws1.Range("F26").Value = 63
$ws1.Range("F27").Value = 1872
$ws1.Range("F29").Value = 1925
$ws1.Range("F31").Value = 1105
$ws1.Range("F32").Value = 268
$ws1.Range("F34").Value = 203
$ws1.Range("F36").Value = 320
$ws1.Range("F37").Value = 63
$ws1.Range("F38").Value = 286
$ws1.Range("F39").Value = 478
$ws1.Range("F40").Value = 5
$ws1.Range("F41").Value = 17
$ws1.Range("F42").Value = 179
$ws1.Range("F43").Value = 1367
$ws1.Range("F44").Value = 291
$ws1.Range("F46").Value = 3
$ws1.Range("F47").Value = 604
$ws1.Range("F49").Value = 294

# Sheet "演出" (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 6

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2657
$ws4.Range("F3").Value = 339
$ws4.Range("F4").Value = 353
$ws4.Range("F5").Value = 1485
$ws4.Range("F7").Value = 1134
$ws4.Range("F8").Value = 326
$ws4.Range("F10").Value = 115
$ws4.Range("F12").Value = 9071
$ws4.Range("F13").Value = 391
$ws4.Range("F14").Value = 2500
$ws4.Range("F16").Value = 2
$ws4.Range("F18").Value = 255
$ws4.Range("F21").Value = 623
$ws4.Range("F22").Value = 1177
$ws4.Range("F24").Value = 2165
$ws4.Range("F25").Value = 1872
$ws4.Range("F26").Value = 1925
$ws4.Range("F28").Value = 1105
$ws4.Range("F29").Value = 268
$ws4.Range("F31").Value = 203
$ws4.Range("F33").Value = 320
$ws4.Range("F34").Value = 63
$ws4.Range("F35").Value = 286
$ws4.Range("F36").Value = 478
$ws4.Range("F37").Value = 17
$ws4.Range("F38").Value = 6
$ws4.Range("F40").Value = 5
$ws4.Range("F41").Value = 17
$ws4.Range("F42").Value = 179
$ws4.Range("F44").Value = 1367
$ws4.Range("F46").Value = 291
$ws4.Range("F48").Value = 604
$ws4.Range("F49").Value = 294
